# Insert a new weekly price record as row 41 ("Feria Lagunitas de Puerto
# Montt" - Albahaca), pushing the existing rows 41-91 down to 42-92.
# This grows the used range from A1:R91 to A1:R92.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 41, shifting rows 41..91 down to 42..92.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new record's values.
$ws.Cells.Item(41, 1).Value  = 4
$ws.Cells.Item(41, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(41, 3).Value  = "Los Lagos"
$ws.Cells.Item(41, 4).Value  = 44579
$ws.Cells.Item(41, 5).Value  = 10
$ws.Cells.Item(41, 6).Value  = 100112052
$ws.Cells.Item(41, 7).Value  = "Albahaca"
$ws.Cells.Item(41, 8).Value  = "Sin especificar"
$ws.Cells.Item(41, 9).Value  = "Primera"
$ws.Cells.Item(41, 10).Value = 120
$ws.Cells.Item(41, 11).Value = 7000
$ws.Cells.Item(41, 12).Value = 7000
$ws.Cells.Item(41, 13).Value = 7000
$ws.Cells.Item(41, 14).Value = "`$/docena de matas"
$ws.Cells.Item(41, 15).Value = "Región Metropolitana"
$ws.Cells.Item(41, 16).Value = 1167
$ws.Cells.Item(41, 17).Value = 6
$ws.Cells.Item(41, 18).Value = "Hortaliza"
